# EPBDS-8161 Create a new property "calculateAllCells" for Spreadsheets
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Properties")

# Insert a new row above the current row 48 ("Concurrent Execution"),
# shifting it (and everything below) down by one row.
$ws.Rows.Item(48).Insert()

# Copy the formatting of the row above (row 47, "Auto Type Discovery") into
# the freshly inserted row so the new row matches the sheet's styling.
$ws.Range("B47:T47").Copy()
$ws.Range("B48:T48").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Rows.Item(48).RowHeight = $ws.Rows.Item(47).RowHeight

# Populate the new property row: "Calculate All Cells" / calculateAllCells
# (cells are written in the same order the new shared strings were added in
# the authoritative edit, so the shared-string table indices line up: first
# the table-type, then the internal name, then the display name, and the
# description last.)
$ws.Range("D48").Value2 = "no"
$ws.Range("E48").Value2 = "Boolean"
$ws.Range("F48").ClearContents()
$ws.Range("G48").Value2 = "no"
$ws.Range("H48").Value2 = "Dev"
$ws.Range("I48").ClearContents()
$ws.Range("J48").ClearContents()
$ws.Range("K48").Value2 = "no"
$ws.Range("L48").Value2 = "no"
$ws.Range("M48").Value2 = "XLS_SPREADSHEET"
$ws.Range("C48").Value2 = "calculateAllCells"
$ws.Range("B48").Value2 = "Calculate All Cells"
$ws.Range("N48").Formula = "=TRUE()"
$ws.Range("O48").ClearContents()
$ws.Range("P48").ClearContents()
$ws.Range("Q48").ClearContents()
$ws.Range("R48").ClearContents()
$ws.Range("S48").Value2 = "MODULE, CATEGORY, TABLE"
$ws.Range("T48").Value2 = "If true calculates all cells in the Spreadsheet, otherwise calculates only cells these are requred for a result. By default = true."

# Reflect the author's final selection on this sheet.
$ws.Activate() | Out-Null
$ws.Range("M46").Select() | Out-Null
